# Applies the edits described by the commit diff to ClassDesignMapping.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "driver" / "log" declaration text in the BaseTest box (N3/N4)
$ws.Range("N3").Value = "WebDriver driver=new Webdriver();"
$ws.Range("N4").Value = "Logger log="

# 2) Un-merge the A13:D13 box so A13 becomes a standalone cell again
$ws.Range("A13:D13").UnMerge()

# 3) Re-arrange the BasePageObject method list in column A (rows 13-16)
$ws.Range("A13").Value = " -void dismissAlert()"
$ws.Range("A14").Value = " -void switchToWindow(String title)"
$ws.Range("A15").Value = ""
$ws.Range("A16").Value = " -void acceptAlert()"

# 4) Move the active selection the way the author left it
$ws.Range("E23").Select()
